$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New parameter rows 78-94 (appended below the existing table) ---
# Shared strings must be created in this exact order (matching the target
# sharedStrings.xml append order): rows 78-89, then 91-94, then 90 last
# (the author typed "Kmc_EC_NCX_C" into A90 after already filling in the
# later rows, so its string landed at the end of the shared-string table).

$ws.Range("A78").Value = "V_a"
$ws.Range("B78").Value = 70

$ws.Range("A79").Value = "V_h"
$ws.Range("B79").Value = -45

$ws.Range("A80").Value = "V_hkinf"
$ws.Range("B80").Value = -40

$ws.Range("A81").Value = "V_m"
$ws.Range("B81").Value = -46

$ws.Range("A82").Value = "V_n"
$ws.Range("B82").Value = -40

$ws.Range("A83").Value = "V_Sinf"
$ws.Range("B83").Value = -78

$ws.Range("A84").Value = "V_tau"
$ws.Range("B84").Value = 90

$ws.Range("A85").Value = "VBar_RyR"
$ws.Range("B85").Value = -20

$ws.Range("A86").Value = "kATP"
$ws.Range("B86").Value = 0.04

$ws.Range("A87").Value = "KmNa_i_NCX"
$ws.Range("B87").Value = 12290

$ws.Range("A88").Value = "KmNa_EC_NCX"
$ws.Range("B88").Value = 87500

$ws.Range("A89").Value = "Kmc_EC_NCX_N"
$ws.Range("B89").Value = 1300

$ws.Range("A91").Value = "g_SLLeak"
$ws.Range("B91").Value = [double]"9.99999999999999912396E-06"

$ws.Range("A92").Value = "L_RyR"
$ws.Range("B92").Value = 500000

$ws.Range("A93").Value = "g0_DHPR"
$ws.Range("B93").Value = [double]"9.39000000000000112355E-02"

$ws.Range("A94").Value = "j0_RyR"
$ws.Range("B94").Value = 43846000

# Row 90 filled in last, so its label string is appended at the very end
# of the shared-string table (index 93) instead of inline (index 88).
$ws.Range("A90").Value = "Kmc_EC_NCX_C"
$ws.Range("B90").Value = 1600

# Row 96: a lone formatted-but-empty cell (matches style index 2 == "0.00E+00").
$ws.Range("D96").NumberFormat = "0.00E+00"

# --- View state: scroll the sheet down and move the selection ---
$ws.Range("A85").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 64
$win.ScrollColumn = 1
